# Scheduled-runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, reflecting newer market-board prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3003
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 4006
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 4006
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -4230

$ws.Range("H14").Value = 3003
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 4006
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 4006
$ws.Range("M14").Value = -1809
$ws.Range("N14").Value = -4388

$ws.Range("H21").Value = 8000
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 9875
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 9875
$ws.Range("M21").Value = -2532
$ws.Range("N21").Value = -10811

$ws.Range("H23").Value = 8000
$ws.Range("I23").Value = 3000
$ws.Range("J23").Value = 9875
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 9875
$ws.Range("M23").Value = -2766
$ws.Range("N23").Value = -10343

$ws.Range("H98").Value = 104167570
$ws.Range("I98").Value = 104167570
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 104167570
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -104166072
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 104167570
$ws.Range("I122").Value = 104167570
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 312502710
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -312500260
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 2098.509
$ws.Range("I137").Value = 2061.7104
$ws.Range("J137").Value = 2180.7646
$ws.Range("K137").Value = 6185.1312
$ws.Range("L137").Value = 6542.293799999999
$ws.Range("M137").Value = -3635.1312
$ws.Range("N137").Value = -11642.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7803212.5
$ws.Range("I32").Value = 2474584.5
$ws.Range("J32").Value = 27785568
$ws.Range("K32").Value = 2474584.5
$ws.Range("L32").Value = 27785568
$ws.Range("M32").Value = -2474297.5
$ws.Range("N32").Value = -27786142

$ws.Range("H74").Value = 1863.7142
$ws.Range("I74").Value = 1720.8125
$ws.Range("J74").Value = 1984.0526
$ws.Range("K74").Value = 1720.8125
$ws.Range("L74").Value = 1984.0526
$ws.Range("M74").Value = -846.8125
$ws.Range("N74").Value = -3732.0526

$ws.Range("H77").Value = 1863.7142
$ws.Range("I77").Value = 1720.8125
$ws.Range("J77").Value = 1984.0526
$ws.Range("K77").Value = 8604.0625
$ws.Range("L77").Value = 9920.262999999999
$ws.Range("M77").Value = -4236.0625
$ws.Range("N77").Value = -18656.263

$ws.Range("H132").Value = 1155496.2
$ws.Range("I132").Value = 2040.5
$ws.Range("J132").Value = 2803290.2
$ws.Range("K132").Value = 6121.5
$ws.Range("L132").Value = 8409870.600000001
$ws.Range("M132").Value = -3591.5
$ws.Range("N132").Value = -8414930.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1999
$ws.Range("I94").Value = 1999
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1548

$ws.Range("H141").Value = 357000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 357000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 357000
$ws.Range("N141").Value = -367360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 457.3125
$ws.Range("I107").Value = 473.16666
$ws.Range("J107").Value = 409.75
$ws.Range("K107").Value = 473.16666
$ws.Range("L107").Value = 409.75
$ws.Range("M107").Value = 1446.83334
$ws.Range("N107").Value = -4249.75

$ws.Range("H132").Value = 11496186
$ws.Range("I132").Value = 1407.238
$ws.Range("J132").Value = 41669984
$ws.Range("K132").Value = 4221.714
$ws.Range("L132").Value = 125009952
$ws.Range("M132").Value = -1691.714
$ws.Range("N132").Value = -125015012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 122245080
$ws.Range("I9").Value = 33333700
$ws.Range("J9").Value = 166700770
$ws.Range("K9").Value = 100001100
$ws.Range("L9").Value = 500102310
$ws.Range("M9").Value = -100000876
$ws.Range("N9").Value = -500102758

$ws.Range("H11").Value = 68.666664
$ws.Range("I11").Value = 62.4
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 187.2
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -47.19999999999999
$ws.Range("N11").Value = -580

$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 6000
$ws.Range("N94").Value = -7352

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1070.4
$ws.Range("I36").Value = 1238
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 1238
$ws.Range("L36").Value = 400
$ws.Range("M36").Value = -753

$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3344

$ws.Range("H97").Value = 2527.5
$ws.Range("I97").Value = 2138.5
$ws.Range("J97").Value = 3500
$ws.Range("K97").Value = 2138.5
$ws.Range("L97").Value = 3500
$ws.Range("M97").Value = -1642.5
$ws.Range("N97").Value = -4492

$ws.Range("H122").Value = 20843730
$ws.Range("I122").Value = 38478372
$ws.Range("J122").Value = 2790.5454
$ws.Range("K122").Value = 115435116
$ws.Range("L122").Value = 8371.636200000001
$ws.Range("M122").Value = -115432666
$ws.Range("N122").Value = -13271.6362

$ws.Range("H126").Value = 2823.7334
$ws.Range("I126").Value = 2172.6667
$ws.Range("J126").Value = 3800.3333
$ws.Range("K126").Value = 6518.000100000001
$ws.Range("L126").Value = 11400.9999
$ws.Range("M126").Value = -4048.000100000001
$ws.Range("N126").Value = -16340.9999

$ws.Range("H132").Value = 5063.2666
$ws.Range("I132").Value = 2546.2903
$ws.Range("J132").Value = 10636.571
$ws.Range("K132").Value = 7638.8709
$ws.Range("L132").Value = 31909.713
$ws.Range("M132").Value = -5108.8709
$ws.Range("N132").Value = -36969.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1675
$ws.Range("I7").Value = 900
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -788

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H22").Value = 4433.3335
$ws.Range("I22").Value = 6500
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 6500
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -6205
$ws.Range("N22").Value = -890

$ws.Range("H27").Value = 4433.3335
$ws.Range("I27").Value = 6500
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 6500
$ws.Range("L27").Value = 300
$ws.Range("M27").Value = -6393
$ws.Range("N27").Value = -514

$ws.Range("H40").Value = 12502395
$ws.Range("I40").Value = 1929.0714
$ws.Range("J40").Value = 41670150
$ws.Range("K40").Value = 1929.0714
$ws.Range("L40").Value = 41670150
$ws.Range("M40").Value = -1793.0714
$ws.Range("N40").Value = -41670422

$ws.Range("H126").Value = 1675
$ws.Range("I126").Value = 900
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 2700
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -230

$ws.Range("H132").Value = 17550320
$ws.Range("I132").Value = 38464264
$ws.Range("J132").Value = 9593.807000000001
$ws.Range("K132").Value = 115392792
$ws.Range("L132").Value = 28781.421
$ws.Range("M132").Value = -115390262
$ws.Range("N132").Value = -33841.421

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 31723.234
$ws.Range("I122").Value = 46900.184
$ws.Range("J122").Value = 3898.8333
$ws.Range("K122").Value = 140700.552
$ws.Range("L122").Value = 11696.4999
$ws.Range("M122").Value = -138250.552

$ws.Range("H126").Value = 2657.1765
$ws.Range("I126").Value = 2380.5833
$ws.Range("J126").Value = 3321
$ws.Range("K126").Value = 7141.749899999999
$ws.Range("L126").Value = 9963
$ws.Range("M126").Value = -4671.749899999999
$ws.Range("N126").Value = -14903

$ws.Range("H136").Value = 8476043
$ws.Range("I136").Value = 13514258
$ws.Range("J136").Value = 2681.7273
$ws.Range("K136").Value = 40542774
$ws.Range("L136").Value = 8045.1819
$ws.Range("M136").Value = -40540224
$ws.Range("N136").Value = -13145.1819
